# [ADDITIONAL SCRAPING] added scraping code for extra browling attributes and excel sheets
#
# 1) "ODI Batting Extra" - drop the trailing empty B:E cells that a row has
#    no real data for (mirrors how the scraper now only emits cells up to
#    the last populated column instead of padding with blank placeholders).
# 2) Add a brand new "ODI Bowling Extra" sheet (MATCH_CODE / MAIDEN_OVERS /
#    PERCENT_WICKETS_OF_ALL) with the scraped bowling-extra data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Clean up "ODI Batting Extra": clear the now-unused trailing cells.
# ---------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")

$cellsToClear = @(
    "B2","C2","D2","E2",
    "E3",
    "B4","C4","D4","E4",
    "B5","C5","D5","E5",
    "B9","C9","D9","E9",
    "B11","C11","D11","E11",
    "B15","C15","D15","E15",
    "B16","C16","D16","E16",
    "C17","D17","E17",
    "B20","C20","D20","E20"
)

foreach ($addr in $cellsToClear) {
    $battingExtra.Range($addr).ClearContents()
}

# ---------------------------------------------------------------------
# 2) Add the new "ODI Bowling Extra" sheet at the end of the workbook.
# ---------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$bowlingExtra = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets.Item($sheetCount))
$bowlingExtra.Name = "ODI Bowling Extra"

# Header row (bold, centered, bordered - matching the other sheets' headers)
$headers = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $bowlingExtra.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Data rows: MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL
# (blank strings are written as blank text cells, not left absent, and every
# value - including ones that look numeric - is stored as text.)
$rows = @(
    @("4351", "0", "20.00%"),
    @("4401", "", ""),
    @("4408", "", ""),
    @("4415", "0", "10.00%"),
    @("4419", "", ""),
    @("4421", "", ""),
    @("4458", "0", "20.00%"),
    @("4459", "0", "20.00%"),
    @("4460", "0", "10.00%"),
    @("4474", "", ""),
    @("4475", "0", "20.00%"),
    @("4478", "", ""),
    @("4487", "0", ""),
    @("4488", "0", ""),
    @("4491", "0", ""),
    @("4524", "", ""),
    @("4526", "0", "10.00%"),
    @("4529", "0", "30.00%"),
    @("4550", "1", "10.00%"),
    @("4619", "", "")
)

$rowIndex = 2
foreach ($rowData in $rows) {
    for ($col = 1; $col -le 3; $col++) {
        $cell = $bowlingExtra.Cells.Item($rowIndex, $col)
        $value = $rowData[$col - 1]
        if ($value -eq "") {
            $cell.Value = ""
        } else {
            $cell.Value = "'" + $value
        }
    }
    $rowIndex++
}
